$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: Apply Garamond font to every run that lives inside a "Normal"
# (i.e. not Heading1/Heading2/Heading3) paragraph. We select each paragraph's
# range but back off the trailing paragraph-mark character so that we do not
# touch the paragraph-mark run properties (w:pPr/w:rPr) - only the actual
# w:r/w:rPr elements should get the w:rFonts addition.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $styleName = $p.Style.NameLocal
    if ($styleName -ne "Heading 1" -and $styleName -ne "Heading 2" -and $styleName -ne "Heading 3") {
        $rng = $p.Range.Duplicate
        if ($rng.End -gt $rng.Start) {
            $rng.MoveEnd(1, -1)
            if ($rng.End -gt $rng.Start) {
                $rng.Font.Name = "Garamond"
            }
        }
    }
}

# ---------------------------------------------------------------------------
# Step 2: Insert a page-break-only paragraph immediately before every
# "Heading 3" paragraph (the "Part N: ..." section headers). We walk from the
# last Heading 3 paragraph to the first so that earlier paragraph indices are
# not disturbed by subsequent insertions.
# ---------------------------------------------------------------------------
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$breakFragment = "<w:p $wNs><w:r><w:br w:type=`"page`"/></w:r></w:p>"

$headingIndexes = New-Object System.Collections.ArrayList
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Style.NameLocal -eq "Heading 3") {
        [void]$headingIndexes.Add($i)
    }
}

for ($j = $headingIndexes.Count - 1; $j -ge 0; $j--) {
    $idx = $headingIndexes[$j]
    $target = $d.Paragraphs.Item($idx)
    $insertPoint = $target.Range.Duplicate
    $insertPoint.Collapse(1)
    $insertPoint.InsertParagraphBefore()
    $newPara = $d.Paragraphs.Item($idx)
    $newPara.Range.InsertXML($breakFragment)
}

Write-Host "Done. Final paragraph count:" $d.Paragraphs.Count
